# Goshop 13st Commit Rev 1.0
# - Update Sheet1 row 16 (2025-03-09 -> 2025-03-12, new amounts)
# - Add sheet "銷售記錄" with the per-file revenue log (incl. the new 03/12 row)
# - Add sheet "銷售總合" with the grand total revenue

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Sheet1: row 16 changes from 2025-03-09 to 2025-03-12 with new values
# ---------------------------------------------------------------------------
$ws1.Range("A16").NumberFormat = "@"
$ws1.Range("A16").Value = "2025-03-12"
$ws1.Range("A16").Style = "Normal"
$ws1.Range("B16").Value = 794.4100000000001
$ws1.Range("C16").Value = 79.60999999999999
$ws1.Range("D16").Value = 714.8000000000001

# ---------------------------------------------------------------------------
# 2. Add "銷售記錄" (sales record) sheet right after Sheet1
# ---------------------------------------------------------------------------
$wsRecord = $wb.Worksheets.Add($null, $ws1)
$wsRecord.Name = "銷售記錄"
# Match Sheet1's outline settings (summary rows below / summary cols right)
$wsRecord.Outline.SummaryRow = 1
$wsRecord.Outline.SummaryColumn = 1
# Match Sheet1's page margins (0.75in/0.75in/1in/1in/0.5in/0.5in - values are
# in points: 1in = 72pt) instead of Excel's 0.7/0.7/0.75/0.75/0.3/0.3 default.
$wsRecord.PageSetup.LeftMargin = 54
$wsRecord.PageSetup.RightMargin = 54
$wsRecord.PageSetup.TopMargin = 72
$wsRecord.PageSetup.BottomMargin = 72
$wsRecord.PageSetup.HeaderMargin = 36
$wsRecord.PageSetup.FooterMargin = 36

$wsRecord.Range("A1").Value = "檔案名"
$wsRecord.Range("B1").Value = "revenue"
# Reuse Sheet1's header style (bold, centered, thin border) instead of
# building a brand-new style entry.
$ws1.Range("A1").Copy()
$wsRecord.Range("A1:B1").PasteSpecial(-4122)

$names = @(
    "goshop_orders_20250214.xlsx",
    "goshop_orders_20250215.xlsx",
    "goshop_orders_20250217.xlsx",
    "goshop_orders_20250218.xlsx",
    "goshop_orders_20250219_liusming@hotmail.com.xlsx",
    "goshop_orders_20250220_liusming@hotmail.com.xlsx",
    "goshop_orders_20250221_liusming@hotmail.com.xlsx",
    "goshop_orders_20250223_liusming@hotmail.com.xlsx",
    "goshop_orders_20250225_liusming@hotmail.com.xlsx",
    "goshop_orders_20250226_liusming@hotmail.com.xlsx",
    "goshop_orders_20250227_liusming@hotmail.com.xlsx",
    "goshop_orders_20250228_liusming@hotmail.com.xlsx",
    "goshop_orders_20250301_liusming@hotmail.com.xlsx",
    "goshop_orders_20250302_liusming@hotmail.com.xlsx",
    "goshop_orders_20250303_liusming@hotmail.com.xlsx",
    "goshop_orders_20250305_liusming@hotmail.com.xlsx",
    "goshop_orders_20250312_liusming@hotmail.com.xlsx"
)
$vals = @(
    412.2200000000001,
    201.18,
    252.46,
    286.16,
    155.76,
    270.26,
    195.84,
    476.4399999999999,
    349.35,
    316.3099999999999,
    319.52,
    145.4399999999999,
    278.08,
    196.36,
    263.1799999999999,
    644.8000000000001,
    714.8000000000001
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $wsRecord.Cells.Item($r, 1).Value = $names[$i]
    $wsRecord.Cells.Item($r, 2).Value = $vals[$i]
}

# ---------------------------------------------------------------------------
# 3. Add "銷售總合" (sales summary) sheet after 銷售記錄
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Add($null, $wsRecord)
$wsSummary.Name = "銷售總合"
$wsSummary.Outline.SummaryRow = 1
$wsSummary.Outline.SummaryColumn = 1
$wsSummary.PageSetup.LeftMargin = 54
$wsSummary.PageSetup.RightMargin = 54
$wsSummary.PageSetup.TopMargin = 72
$wsSummary.PageSetup.BottomMargin = 72
$wsSummary.PageSetup.HeaderMargin = 36
$wsSummary.PageSetup.FooterMargin = 36

$wsSummary.Range("A1").Value = "總收入"
$ws1.Range("A1").Copy()
$wsSummary.Range("A1").PasteSpecial(-4122)

$wsSummary.Range("A2").Value = 5478.16

# ---------------------------------------------------------------------------
# Keep Sheet1 as the active sheet (matches the original workbook's activeTab)
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
